# Adds a new employee record (row 25 / "nicolas") to the bottom of the
# Hoja1 table, matching the row immediately above it for number formatting
# (currency-style integer on salary, short-date on ingreso), and moves the
# active selection to J17 (as left by the author after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$newRow = 21

$ws.Range("A$newRow").Value = 25
$ws.Range("B$newRow").Value = "nicolas"
$ws.Range("C$newRow").Value = 4354653
$ws.Range("D$newRow").Value = "Heavens Fruits SAS"
$ws.Range("E$newRow").Value = "Activo"
$ws.Range("F$newRow").Value = "Gestion TI"
$ws.Range("G$newRow").Value = 3
$ws.Range("H$newRow").Value = 3100000
$ws.Range("I$newRow").Value = "No"
$ws.Range("J$newRow").Value = 44991

# Match the number formatting (salary / ingreso) of the row above by
# copying its formatting only, so we reuse the existing cell styles
# instead of minting new ones.
$ws.Range("H20").Copy()
$ws.Range("H$newRow").PasteSpecial(-4122)

$ws.Range("J20").Copy()
$ws.Range("J$newRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Leave the selection where the author left it after entering the row.
$ws.Range("J17").Select() | Out-Null
